$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 8000
$ws.Range("I34").Value = 8000
$ws.Range("K34").Value = 8000
$ws.Range("M34").Value = -7797
$ws.Range("H36").Value = 8000
$ws.Range("I36").Value = 8000
$ws.Range("K36").Value = 8000
$ws.Range("M36").Value = -7285
$ws.Range("H43").Value = 3019.6
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 3019.6
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 3019.6
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -3157.6
$ws.Range("H62").Value = 5391.5
$ws.Range("I62").Value = 5474.75
$ws.Range("K62").Value = 5474.75
$ws.Range("M62").Value = -4850.75
$ws.Range("H65").Value = 5391.5
$ws.Range("I65").Value = 5474.75
$ws.Range("K65").Value = 27373.75
$ws.Range("M65").Value = -24253.75
$ws.Range("H76").Value = 8921.571
$ws.Range("I76").Value = 10243.667
$ws.Range("K76").Value = 10243.667
$ws.Range("M76").Value = -9928.666999999999
$ws.Range("H79").Value = 8921.571
$ws.Range("I79").Value = 10243.667
$ws.Range("K79").Value = 10243.667
$ws.Range("M79").Value = -9151.666999999999
$ws.Range("H80").Value = 2802489
$ws.Range("J80").Value = 6803590
$ws.Range("L80").Value = 20410770
$ws.Range("N80").Value = -20412766
$ws.Range("H83").Value = 2802489
$ws.Range("J83").Value = 6803590
$ws.Range("L83").Value = 61232310
$ws.Range("N83").Value = -61242294
$ws.Range("H86").Value = 3186.9285
$ws.Range("I86").Value = 2312.5557
$ws.Range("K86").Value = 2312.5557
$ws.Range("M86").Value = -1189.5557
$ws.Range("H89").Value = 3186.9285
$ws.Range("I89").Value = 2312.5557
$ws.Range("K89").Value = 11562.7785
$ws.Range("M89").Value = -5946.7785
$ws.Range("H92").Value = 2018.625
$ws.Range("I92").Value = 1539.091
$ws.Range("J92").Value = 3073.6
$ws.Range("K92").Value = 1539.091
$ws.Range("L92").Value = 3073.6
$ws.Range("M92").Value = -291.0909999999999
$ws.Range("N92").Value = -5569.6
$ws.Range("H98").Value = 15625998
$ws.Range("I98").Value = 17857884
$ws.Range("K98").Value = 17857884
$ws.Range("M98").Value = -17856386
$ws.Range("H116").Value = 18208
$ws.Range("I116").Value = 9332.666999999999
$ws.Range("J116").Value = 21166.445
$ws.Range("K116").Value = 9332.666999999999
$ws.Range("L116").Value = 21166.445
$ws.Range("M116").Value = -5890.666999999999
$ws.Range("N116").Value = -28050.445
$ws.Range("H122").Value = 15625998
$ws.Range("I122").Value = 17857884
$ws.Range("K122").Value = 53573652
$ws.Range("M122").Value = -53571202
$ws.Range("H137").Value = 2220.5454
$ws.Range("I137").Value = 1491.091
$ws.Range("J137").Value = 3679.4546
$ws.Range("K137").Value = 4473.272999999999
$ws.Range("L137").Value = 11038.3638
$ws.Range("M137").Value = -1923.272999999999
$ws.Range("N137").Value = -16138.3638

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 556.5
$ws.Range("I5").Value = 556.5
$ws.Range("K5").Value = 556.5
$ws.Range("M5").Value = -444.5
$ws.Range("H32").Value = 6885.353
$ws.Range("I32").Value = 6203.271
$ws.Range("J32").Value = 17798.666
$ws.Range("K32").Value = 6203.271
$ws.Range("L32").Value = 17798.666
$ws.Range("M32").Value = -5916.271
$ws.Range("N32").Value = -18372.666
$ws.Range("H60").Value = 89397.664
$ws.Range("I60").Value = 89397.664
$ws.Range("K60").Value = 89397.664
$ws.Range("M60").Value = -88664.664
$ws.Range("H110").Value = 2757.9023
$ws.Range("I110").Value = 2606.1428
$ws.Range("J110").Value = 3084.7693
$ws.Range("K110").Value = 2606.1428
$ws.Range("L110").Value = 3084.7693
$ws.Range("M110").Value = -561.1428000000001
$ws.Range("N110").Value = -7174.7693

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 556.5
$ws.Range("I4").Value = 556.5
$ws.Range("K4").Value = 556.5
$ws.Range("M4").Value = -441.5
$ws.Range("H20").Value = 2834.5642
$ws.Range("I20").Value = 2856.423
$ws.Range("K20").Value = 2856.423
$ws.Range("M20").Value = -2609.423
$ws.Range("H94").Value = 2803.5557
$ws.Range("I94").Value = 3538.8333
$ws.Range("J94").Value = 1333
$ws.Range("K94").Value = 3538.8333
$ws.Range("L94").Value = 1333
$ws.Range("M94").Value = -3087.8333
$ws.Range("N94").Value = -2235

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 27048884
$ws.Range("I31").Value = 30326782
$ws.Range("J31").Value = 6220
$ws.Range("K31").Value = 30326782
$ws.Range("L31").Value = 6220
$ws.Range("M31").Value = -30326487
$ws.Range("N31").Value = -6810
$ws.Range("H34").Value = 27048884
$ws.Range("I34").Value = 30326782
$ws.Range("J34").Value = 6220
$ws.Range("K34").Value = 30326782
$ws.Range("L34").Value = 6220
$ws.Range("M34").Value = -30326580
$ws.Range("N34").Value = -6624
$ws.Range("H41").Value = 35000
$ws.Range("I41").Value = 35000
$ws.Range("K41").Value = 35000
$ws.Range("M41").Value = -34572
$ws.Range("H59").Value = 107499.5
$ws.Range("I59").Value = 0
$ws.Range("K59").Value = 0
$ws.Range("M59").ClearContents()
$ws.Range("H62").Value = 13899926
$ws.Range("I62").Value = 7410.8335
$ws.Range("J62").Value = 27792440
$ws.Range("K62").Value = 7410.8335
$ws.Range("L62").Value = 27792440
$ws.Range("M62").Value = -6786.8335
$ws.Range("N62").Value = -27793688
$ws.Range("H65").Value = 13899926
$ws.Range("I65").Value = 7410.8335
$ws.Range("J65").Value = 27792440
$ws.Range("K65").Value = 37054.1675
$ws.Range("L65").Value = 138962200
$ws.Range("M65").Value = -33934.1675
$ws.Range("N65").Value = -138968440
$ws.Range("H122").Value = 2420
$ws.Range("I122").Value = 2317.1667
$ws.Range("K122").Value = 6951.500100000001
$ws.Range("M122").Value = -4501.500100000001
$ws.Range("H141").Value = 1000000
$ws.Range("J141").Value = 1000000
$ws.Range("L141").Value = 1000000
$ws.Range("N141").Value = -1010360

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 78.57143000000001
$ws.Range("H59").Value = 17692.666
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 17692.666
$ws.Range("K59").Value = 0
$ws.Range("L59").Value = 53077.99800000001
$ws.Range("M59").ClearContents()
$ws.Range("N59").Value = -54157.99800000001
$ws.Range("H69").Value = 12475
$ws.Range("I69").Value = 3294.75
$ws.Range("J69").Value = 21655.25
$ws.Range("K69").Value = 9884.25
$ws.Range("L69").Value = 64965.75
$ws.Range("M69").Value = -9073.25
$ws.Range("N69").Value = -66587.75
$ws.Range("H72").Value = 12475
$ws.Range("I72").Value = 3294.75
$ws.Range("J72").Value = 21655.25
$ws.Range("K72").Value = 29652.75
$ws.Range("L72").Value = 194897.25
$ws.Range("M72").Value = -25596.75
$ws.Range("N72").Value = -203009.25

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H45").Value = 85662.75
$ws.Range("J45").Value = 85662.75
$ws.Range("L45").Value = 85662.75
$ws.Range("N45").Value = -86780.75
$ws.Range("H97").Value = 2153.1052
$ws.Range("I97").Value = 1999.3334
$ws.Range("J97").Value = 2416.7144
$ws.Range("K97").Value = 1999.3334
$ws.Range("L97").Value = 2416.7144
$ws.Range("M97").Value = -1503.3334
$ws.Range("N97").Value = -3408.7144
$ws.Range("H113").Value = 3707301.5
$ws.Range("I113").Value = 4336.6665
$ws.Range("J113").Value = 9261749
$ws.Range("K113").Value = 4336.6665
$ws.Range("L113").Value = 9261749
$ws.Range("M113").Value = -2166.6665
$ws.Range("N113").Value = -9266089
$ws.Range("H122").Value = 5509.8965
$ws.Range("I122").Value = 4975.1763
$ws.Range("K122").Value = 14925.5289
$ws.Range("M122").Value = -12475.5289

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6464
$ws.Range("I22").Value = 13276.777
$ws.Range("K22").Value = 13276.777
$ws.Range("M22").Value = -12981.777
$ws.Range("H27").Value = 6464
$ws.Range("I27").Value = 13276.777
$ws.Range("K27").Value = 13276.777
$ws.Range("M27").Value = -13169.777
$ws.Range("H55").Value = 917.9
$ws.Range("I55").Value = 376.57144
$ws.Range("K55").Value = 376.57144
$ws.Range("M55").Value = -203.57144
$ws.Range("H132").Value = 3936.95
$ws.Range("I132").Value = 2302.1304
$ws.Range("K132").Value = 6906.3912
$ws.Range("M132").Value = -4376.3912

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1389.3572
$ws.Range("I81").Value = 1856.4286
$ws.Range("K81").Value = 3712.8572
$ws.Range("M81").Value = -2651.8572
$ws.Range("H84").Value = 1389.3572
$ws.Range("I84").Value = 1856.4286
$ws.Range("K84").Value = 18564.286
$ws.Range("M84").Value = -13260.286
$ws.Range("H107").Value = 4224.237
$ws.Range("I107").Value = 2233.0454
$ws.Range("J107").Value = 6962.125
$ws.Range("K107").Value = 6699.1362
$ws.Range("L107").Value = 20886.375
$ws.Range("M107").Value = -4779.1362
$ws.Range("N107").Value = -24726.375
$ws.Range("H122").Value = 1905.2727
$ws.Range("I122").Value = 1511.56
$ws.Range("J122").Value = 3135.625
$ws.Range("K122").Value = 4534.68
$ws.Range("L122").Value = 9406.875
$ws.Range("M122").Value = -2084.68
$ws.Range("N122").Value = -14306.875
